$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 'I1rqmbET'
$ws.Range("B5").Value = '31/10/2024'
$ws.Range("C5").Value = '11:50'
$ws.Range("D5").Value = 'SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE'
$ws.Range("E5").Value = 'Al Shabab'
$ws.Range("F5").Value = 'Al Wehda'
$ws.Range("G5").Value = 1.42
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 6.5
$ws.Range("J5").Value = 1.91
$ws.Range("K5").Value = 2.38
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 1.03
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 4.33
$ws.Range("Q5").Value = 1.67
$ws.Range("R5").Value = 2.15
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 3.25
$ws.Range("U5").Value = 1.91
$ws.Range("V5").Value = 1.8
$ws.Range("W5").Value = 7.5
$ws.Range("X5").Value = 7
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 9.5
$ws.Range("AA5").Value = 12
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 13
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 51
$ws.Range("AG5").Value = 700
$ws.Range("AH5").Value = 17
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 21
$ws.Range("AK5").Value = 67
$ws.Range("AL5").Value = 51
$ws.Range("AM5").Value = 51
$ws.Range("AN5").Value = 3.4
$ws.Range("AO5").Value = 7
$ws.Range("AP5").Value = 19
$ws.Range("AQ5").Value = 19
$ws.Range("AR5").Value = 41
$ws.Range("AS5").Value = 126
$ws.Range("AT5").Value = 3.25
$ws.Range("AU5").Value = 9
$ws.Range("AV5").Value = 51
$ws.Range("AW5").Value = 8
$ws.Range("AX5").Value = 34
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 126
$ws.Range("BA5").Value = 151
$ws.Range("BB5").Value = 500
$ws.Range("BC5").Value = 81
$ws.Range("BD5").Value = 81
